# Capacity Supply Curve.xlsx update
# - About sheet: bump the "Date Created" stamp
# - CSC-CSCSoCECBiaSY sheet: set "Share of Cost Effective Capacity Built in a
#   Single Year" to 100% (1) for every resource except crude oil / heavy-or-
#   residual-fuel-oil (rows 16-17), which stay at 0, across all year columns
#   (B:AE)
# - Re-point the active sheet's selection/scroll to the bottom of the table

$wb = $excel.ActiveWorkbook

# --- About sheet: update "Date Created" (C1) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = "03/11/2024"

# --- CSC-CSCSoCECBiaSY: bump share-of-cost-effective-capacity to 1 ---
$wsCSC = $wb.Worksheets.Item("CSC-CSCSoCECBiaSY")
$wsCSC.Range("B2:AE15").Value = 1
$wsCSC.Range("B18:AE25").Value = 1

# --- Update the active sheet's selection to match the edited region ---
$wsCSC.Activate()
$wsCSC.Range("B18:AE25").Select() | Out-Null
